# Insert a new data row at row 170 (pushing the existing rows 170-301 down
# to 171-302) and populate it with a new "Brócoli" price observation for
# Macroferia Regional de Talca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 170..301 down to 171..302, leaving a blank row 170.
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new record.
$ws.Cells.Item(170, 1).Value = 5
$ws.Cells.Item(170, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(170, 3).Value = "Maule"
$ws.Cells.Item(170, 4).Value = 44651
$ws.Cells.Item(170, 5).Value = 7
$ws.Cells.Item(170, 6).Value = 100112023
$ws.Cells.Item(170, 7).Value = "Brócoli"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 2000
$ws.Cells.Item(170, 11).Value = 700
$ws.Cells.Item(170, 12).Value = 700
$ws.Cells.Item(170, 13).Value = 700
$ws.Cells.Item(170, 14).Value = "`$/unidad"
$ws.Cells.Item(170, 15).Value = "Región del Maule"
$ws.Cells.Item(170, 16).Value = 700
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = "Hortaliza"
